$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the date-column (A) number formatting down to the new rows ---
# Column A uses a date display format (style applied to A2:A547); reuse it
# for the newly appended rows instead of building a fresh number format.
$ws.Range("A547").Copy()
$ws.Range("A548:A554").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 548 ---
$ws.Range("A548").Value = 45865
$ws.Range("B548").Value = "Flowering"
$ws.Range("C548").Value = "Large"
$ws.Range("D548").Value = 70
$ws.Range("E548").Value = 82
$ws.Range("F548").Formula = "=ABS(D548-E548)"
$ws.Range("G548").Value = 0.25
$ws.Range("H548").Value = 0.1
$ws.Range("I548").Value = "No"
$ws.Range("J548").Value = 2
$ws.Range("K548").Value = "Bright"
$ws.Range("L548").Value = 4
$ws.Range("M548").Value = 0.82
$ws.Range("N548").Value = 75
$ws.Range("O548").Value = 30.06
$ws.Range("P548").Value = 15
$ws.Range("Q548").Value = 0.74
$ws.Range("R548").Value = 6.2
$ws.Range("S548").Value = 44
$ws.Range("T548").Value = 26

# --- Row 549 ---
$ws.Range("A549").Value = 45865
$ws.Range("B549").Value = "Nonflowering"
$ws.Range("C549").Value = "Medium"
$ws.Range("D549").Value = 70
$ws.Range("E549").Value = 82
$ws.Range("F549").Formula = "=ABS(D549-E549)"
$ws.Range("G549").Value = 0.25
$ws.Range("H549").Value = 0.05
$ws.Range("I549").Value = "No"
$ws.Range("J549").Value = 3
$ws.Range("K549").Value = "Bright"
$ws.Range("L549").Value = 4
$ws.Range("M549").Value = 0.82
$ws.Range("N549").Value = 75
$ws.Range("O549").Value = 30.06
$ws.Range("P549").Value = 15
$ws.Range("Q549").Value = 0.74
$ws.Range("R549").Value = 6.2
$ws.Range("S549").Value = 44
$ws.Range("T549").Value = 26

# --- Row 550 ---
$ws.Range("A550").Value = 45865
$ws.Range("B550").Value = "Nonflowering"
$ws.Range("C550").Value = "Small"
$ws.Range("D550").Value = 70
$ws.Range("E550").Value = 82
$ws.Range("F550").Formula = "=ABS(D550-E550)"
$ws.Range("G550").Value = 0.25
$ws.Range("H550").Value = 0.1
$ws.Range("I550").Value = "No"
$ws.Range("J550").Value = 3
$ws.Range("K550").Value = "Neutral"
$ws.Range("L550").Value = 4
$ws.Range("M550").Value = 0.82
$ws.Range("N550").Value = 75
$ws.Range("O550").Value = 30.06
$ws.Range("P550").Value = 15
$ws.Range("Q550").Value = 0.74
$ws.Range("R550").Value = 6.2
$ws.Range("S550").Value = 44
$ws.Range("T550").Value = 26

# --- Row 551 ---
$ws.Range("A551").Value = 45865
$ws.Range("B551").Value = "Nonflowering"
$ws.Range("C551").Value = "Medium"
$ws.Range("D551").Value = 70
$ws.Range("E551").Value = 82
$ws.Range("F551").Formula = "=ABS(D551-E551)"
$ws.Range("G551").Value = 0.25
$ws.Range("H551").Value = 0.2
$ws.Range("I551").Value = "No"
$ws.Range("J551").Value = 3
$ws.Range("K551").Value = "Neutral"
$ws.Range("L551").Value = 4
$ws.Range("M551").Value = 0.82
$ws.Range("N551").Value = 75
$ws.Range("O551").Value = 30.06
$ws.Range("P551").Value = 15
$ws.Range("Q551").Value = 0.74
$ws.Range("R551").Value = 6.2
$ws.Range("S551").Value = 44
$ws.Range("T551").Value = 26

# --- Row 552 ---
$ws.Range("A552").Value = 45865
$ws.Range("B552").Value = "Nonflowering"
$ws.Range("C552").Value = "Medium"
$ws.Range("D552").Value = 70
$ws.Range("E552").Value = 82
$ws.Range("F552").Formula = "=ABS(D552-E552)"
$ws.Range("G552").Value = 0.25
$ws.Range("H552").Value = 0.2
$ws.Range("I552").Value = "No"
$ws.Range("J552").Value = 3
$ws.Range("K552").Value = "Neutral"
$ws.Range("L552").Value = 4
$ws.Range("M552").Value = 0.82
$ws.Range("N552").Value = 75
$ws.Range("O552").Value = 30.06
$ws.Range("P552").Value = 15
$ws.Range("Q552").Value = 0.74
$ws.Range("R552").Value = 6.2
$ws.Range("S552").Value = 44
$ws.Range("T552").Value = 26

# --- Row 553 ---
$ws.Range("A553").Value = 45865
$ws.Range("B553").Value = "Nonflowering"
$ws.Range("C553").Value = "Large"
$ws.Range("D553").Value = 70
$ws.Range("E553").Value = 82
$ws.Range("F553").Formula = "=ABS(D553-E553)"
$ws.Range("G553").Value = 0.25
$ws.Range("H553").Value = 0.15
$ws.Range("I553").Value = "No"
$ws.Range("J553").Value = 4
$ws.Range("K553").Value = "Bright"
$ws.Range("L553").Value = 4
$ws.Range("M553").Value = 0.82
$ws.Range("N553").Value = 75
$ws.Range("O553").Value = 30.06
$ws.Range("P553").Value = 15
$ws.Range("Q553").Value = 0.74
$ws.Range("R553").Value = 6.2
$ws.Range("S553").Value = 44
$ws.Range("T553").Value = 26

# --- Row 554 ---
$ws.Range("A554").Value = 45865
$ws.Range("B554").Value = "Tree"
$ws.Range("C554").Value = "Medium"
$ws.Range("D554").Value = 70
$ws.Range("E554").Value = 82
$ws.Range("F554").Formula = "=ABS(D554-E554)"
$ws.Range("G554").Value = 0.25
$ws.Range("H554").Value = 0.55000000000000004
$ws.Range("I554").Value = "No"
$ws.Range("J554").Value = 1
$ws.Range("K554").Value = "Bright"
$ws.Range("L554").Value = 4
$ws.Range("M554").Value = 0.82
$ws.Range("N554").Value = 75
$ws.Range("O554").Value = 30.06
$ws.Range("P554").Value = 15
$ws.Range("Q554").Value = 0.74
$ws.Range("R554").Value = 6.2
$ws.Range("S554").Value = 44
$ws.Range("T554").Value = 26

# --- Update the saved view state (scroll position + active selection) ---
$null = $ws.Range("I555").Select()
$excel.ActiveWindow.ScrollRow = 530
